# Update "paises" COVID-19 data: refresh timestamp, refresh per-country
# statistics, and update country labels where refreshed ranking caused
# two (or more) adjacent rows to swap positions in the (descending by
# total cases) table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 30 de Mayo de 2020 a las 13:10"

$ws.Range("B12").Value = 174355
$ws.Range("C12").Value = 864
$ws.Range("D12").Value = 82804
$ws.Range("E12").Value = 86569
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 4982

$ws.Range("A14").Value = "Iran"
$ws.Range("B14").Value = 148950
$ws.Range("C14").Value = 2282
$ws.Range("D14").Value = 116827
$ws.Range("E14").Value = 24389
$ws.Range("G14").Value = 57
$ws.Range("H14").Value = 7734

$ws.Range("A15").Value = "Peru"
$ws.Range("B15").Value = 148285
$ws.Range("D15").Value = 62791
$ws.Range("E15").Value = 81264
$ws.Range("H15").Value = 4230

$ws.Range("B32").Value = 30845
$ws.Range("C32").Value = 17
$ws.Range("D32").Value = 28400
$ws.Range("E32").Value = 526

$ws.Range("B41").Value = 19133
$ws.Range("C41").Value = 151
$ws.Range("D41").Value = 13046
$ws.Range("E41").Value = 4834

$ws.Range("B53").Value = 10740
$ws.Range("C53").Value = 291
$ws.Range("D53").Value = 5811
$ws.Range("E53").Value = 4914

$ws.Range("B64").Value = 7740
$ws.Range("C64").Value = 26
$ws.Range("D64").Value = 5377
$ws.Range("E64").Value = 2160
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 203

$ws.Range("A78").Value = "Senegal"
$ws.Range("B78").Value = 3535
$ws.Range("C78").Value = 106
$ws.Range("D78").Value = 1761
$ws.Range("E78").Value = 1732
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 42

$ws.Range("A79").Value = "Uzbekistan"
$ws.Range("B79").Value = 3513
$ws.Range("C79").Value = 45
$ws.Range("D79").Value = 2728
$ws.Range("E79").Value = 771
$ws.Range("H79").Value = 14

$ws.Range("B87").Value = 2494
$ws.Range("C87").Value = 9
$ws.Range("D87").Value = 1831
$ws.Range("E87").Value = 510

$ws.Range("B101").Value = 1561
$ws.Range("C101").Value = 3
$ws.Range("E101").Value = 770

$ws.Range("A105").Value = "Nepal"
$ws.Range("B105").Value = 1401
$ws.Range("C105").Value = 189
$ws.Range("D105").Value = 219
$ws.Range("E105").Value = 1176
$ws.Range("H105").Value = 6

$ws.Range("A106").Value = "Venezuela"
$ws.Range("B106").Value = 1370
$ws.Range("D106").Value = 302
$ws.Range("E106").Value = 1054
$ws.Range("H106").Value = 14

$ws.Range("A107").Value = "Guinea Ecuatorial"
$ws.Range("B107").Value = 1306
$ws.Range("D107").Value = 200
$ws.Range("E107").Value = 1094
$ws.Range("H107").Value = 12

$ws.Range("A108").Value = "Guinea-Bisau"
$ws.Range("B108").Value = 1256
$ws.Range("D108").Value = 42
$ws.Range("E108").Value = 1206
$ws.Range("H108").Value = 8

$ws.Range("A109").Value = "Mali"
$ws.Range("B109").Value = 1226
$ws.Range("D109").Value = 669
$ws.Range("E109").Value = 484
$ws.Range("H109").Value = 73

$ws.Range("A115").Value = "Etiopia"
$ws.Range("B115").Value = 1063
$ws.Range("C115").Value = 95
$ws.Range("D115").Value = 208
$ws.Range("E115").Value = 847
$ws.Range("H115").Value = 8

$ws.Range("A116").Value = "Zambia"
$ws.Range("B116").Value = 1057
$ws.Range("D116").Value = 779
$ws.Range("E116").Value = 271
$ws.Range("H116").Value = 7

$ws.Range("A117").Value = "Costa Rica"
$ws.Range("B117").Value = 1022
$ws.Range("D117").Value = 653
$ws.Range("E117").Value = 359

$ws.Range("A118").Value = "Sudan del Sur"
$ws.Range("B118").Value = 994
$ws.Range("D118").Value = 6
$ws.Range("E118").Value = 978
$ws.Range("H118").Value = 10

$ws.Range("A129").Value = "Madagascar"
$ws.Range("B129").Value = 758
$ws.Range("C129").Value = 60
$ws.Range("D129").Value = 165
$ws.Range("E129").Value = 587
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = 6

$ws.Range("A130").Value = "Georgia"
$ws.Range("B130").Value = 757
$ws.Range("C130").Value = 11
$ws.Range("D130").Value = 600
$ws.Range("E130").Value = 145
$ws.Range("H130").Value = 12

$ws.Range("A131").Value = "Jordania"
$ws.Range("B131").Value = 730
$ws.Range("D131").Value = 507
$ws.Range("E131").Value = 214
$ws.Range("H131").Value = 9

$ws.Range("A132").Value = "Crucero"
$ws.Range("B132").Value = 712
$ws.Range("D132").Value = 651
$ws.Range("E132").Value = 48
$ws.Range("H132").Value = 13

$ws.Range("D133").Value = 357
$ws.Range("E133").Value = 272

$ws.Range("B134").Value = 618
$ws.Range("C134").Value = 2
$ws.Range("D134").Value = 525
$ws.Range("E134").Value = 84

$ws.Range("A164").Value = "Zimbabue"
$ws.Range("B164").Value = 160
$ws.Range("C164").Value = 11
$ws.Range("D164").Value = 29
$ws.Range("E164").Value = 127
$ws.Range("H164").Value = 4

$ws.Range("A165").Value = "Guyana"
$ws.Range("B165").Value = 150
$ws.Range("D165").Value = 67
$ws.Range("E165").Value = 72
$ws.Range("H165").Value = 11

$ws.Range("A198").Value = "Fiyi"
$ws.Range("D198").Value = 15
$ws.Range("H198").Value = 0

$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 14
$ws.Range("H199").Value = 1

$ws.Range("A210").Value = "Seychelles"
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = "Montserrat"
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("D213").Value = 8
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("D214").Value = 7
$ws.Range("H214").Value = 1
